$wb = $excel.ActiveWorkbook

# --- Hoja2: update matrix values ------------------------------------------------
$ws2 = $wb.Worksheets.Item("Hoja2")

$ws2.Range("G5").Value = 15
$ws2.Range("E7").Value = 15
$ws2.Range("H8").Value = 10

# --- Hoja4: duplicate of Hoja2, inserted right after it (before PageSetup is
# applied to Hoja2, so the new sheet does not inherit it) -----------------------
$ws2.Copy([System.Reflection.Missing]::Value, $ws2)
$ws4 = $wb.Worksheets.Item("Hoja2 (2)")
$ws4.Name = "Hoja4"

# The new sheet differs from Hoja2 by a single cell
$ws4.Range("H8").Value = 5

# Page setup (A4 / portrait) gets recorded on Hoja2 only, after the copy
$ps2 = $ws2.PageSetup
$ps2.PaperSize = 9
$ps2.Orientation = 1

# --- Selections / active views ---------------------------------------------------
$ws4.Activate()
$ws4.Range("E12").Select()

$ws3 = $wb.Worksheets.Item("Hoja3")
$ws3.Activate()
$ws3.Range("F21").Select()

# Restore Hoja2 as the active sheet/tab, matching the saved selection on it
$ws2.Activate()
$ws2.Range("G5").Select()
